$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Formatting: bold font, thin box border, centered horizontal / top vertical alignment
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop
$b1.Borders.LineStyle = 1        # xlContinuous
$b1.Borders.Weight = 2           # xlThin

# Apply the identical formatting to A2 by copying it, so both cells share one style
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
